$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank out the data rows (2-4), replacing values with a single space,
# keeping the header row (row 1) intact.
$ws.Range("A2:G4").Value = " "

# Remove row 5 entirely (shifts dimension from A1:G5 to A1:G4).
$ws.Rows.Item(5).Delete()

# Adjust column widths per the updated layout.
# (Excel's ColumnWidth uses character units that serialize to the sheet's
#  stored <col width> plus a fixed ~0.8333 padding offset, so subtract
#  that offset here to land on the exact target stored widths: 15/12/13.)
$ws.Columns.Item(1).ColumnWidth = 15 - 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 12 - 0.8333333333333334
$ws.Columns.Item(5).ColumnWidth = 13 - 0.8333333333333334
